$wb = $excel.ActiveWorkbook

# GLOBAL RESULTS
$ws = $wb.Worksheets.Item("GLOBAL RESULTS")
$ws.Range("C2").Value = 1375.293247
$ws.Range("C4").Value = 106.0
$ws.Range("C6").Value = 27020.304767035734
$ws.Range("C7").Value = 26596.304767035734
$ws.Range("C8").Value = 26209.695624024658
$ws.Range("C9").Value = 7632.0
$ws.Range("C11").Value = 5708.403422177888
$ws.Range("C13").Value = 21311.901344857844
$ws.Range("C14").Value = 20887.901344857844
$ws.Range("C15").Value = 13679.901344857848
$ws.Range("C16").Value = 13450.357699357846
$ws.Range("C17").Value = 12829.93369935785
$ws.Range("C21").Value = 264978.6717436509
$ws.Range("C22").Value = 260820.65214365092
$ws.Range("C23").Value = 257029.31159134133
$ws.Range("C24").Value = 74844.35279999998
$ws.Range("C27").Value = 208998.3573235501
$ws.Range("C28").Value = 204840.3377235501
$ws.Range("C29").Value = 134154.00452355016
$ws.Range("C30").Value = 131902.95033240755
$ws.Range("C31").Value = 125818.66931280764

# FUSELAGE
$ws = $wb.Worksheets.Item("FUSELAGE")
$ws.Range("C2").Value = 6031.124999999998
$ws.Range("D6").Value = -44.09003295405085
$ws.Range("C7").Value = 2478.0
$ws.Range("D7").Value = -58.91313809612633
$ws.Range("C8").Value = 3140.0
$ws.Range("D8").Value = -47.93674480300108
$ws.Range("C9").Value = 2968.0
$ws.Range("D9").Value = -50.78861738067109
$ws.Range("D10").Value = -58.69758958734895
$ws.Range("D11").Value = -38.684739580094906
$ws.Range("C12").Value = 3024.5
$ws.Range("D12").Value = -49.85181040021553

# WING
$ws = $wb.Worksheets.Item("WING")
$ws.Range("C2").Value = 4020.749999999999
$ws.Range("C7").Value = 2488.0
$ws.Range("D7").Value = -38.12099732636945
$ws.Range("C8").Value = 2738.0
$ws.Range("D8").Value = -31.903251880867984
$ws.Range("C9").Value = 3913.0
$ws.Range("D9").Value = -2.6798482870111076
$ws.Range("D10").Value = -42.82161288316855
$ws.Range("C11").Value = 4931.0
$ws.Range("D11").Value = 22.63881116707085
$ws.Range("C12").Value = 3506.0
$ws.Range("D12").Value = -12.80233787228749
$ws.Range("C13").Value = 2839.2857142857138
$ws.Range("D13").Value = -29.384176726090526

# HORIZONTAL TAIL
$ws = $wb.Worksheets.Item("HORIZONTAL TAIL")
$ws.Range("C2").Value = 603.1124999999997
$ws.Range("D7").Value = -57.22191133494992
$ws.Range("D8").Value = -75.29482476320752
$ws.Range("C9").Value = 140.0
$ws.Range("D9").Value = -76.78708367012787
$ws.Range("C10").Value = 182.33333333333331
$ws.Range("D10").Value = -69.76793992276176

# VERTICAL TAIL
$ws = $wb.Worksheets.Item("VERTICAL TAIL")
$ws.Range("C2").Value = 603.1124999999997
$ws.Range("D7").Value = -45.283840079587115
$ws.Range("C8").Value = 262.0
$ws.Range("D8").Value = -56.55868515409644
$ws.Range("C9").Value = 296.0
$ws.Range("D9").Value = -50.92126261684178

# NACELLES
$ws = $wb.Worksheets.Item("NACELLES")
$ws.Range("C2").Value = 603.1124999999997
$ws.Range("C3").Value = 514.6666666666665
$ws.Range("D3").Value = -71.55496602117255
$ws.Range("D9").Value = -82.75611929780926
$ws.Range("D10").Value = 75.09171174532123
$ws.Range("D11").Value = -36.330286638065004
$ws.Range("C12").Value = 257.33333333333326
$ws.Range("D16").Value = -82.75611929780926
$ws.Range("D17").Value = 75.09171174532123
$ws.Range("D18").Value = -36.330286638065004
$ws.Range("C19").Value = 257.33333333333326

# LANDING GEARS
$ws = $wb.Worksheets.Item("LANDING GEARS")
$ws.Range("C2").Value = 1608.2999999999997
$ws.Range("C5").Value = 874.0
$ws.Range("D5").Value = -45.65690480631722
$ws.Range("C6").Value = 1081.0
$ws.Range("D6").Value = -32.78617173412919
$ws.Range("C7").Value = 1221.0
$ws.Range("D7").Value = -24.081328110427144
$ws.Range("C8").Value = 1095.0
$ws.Range("D8").Value = -31.915687371758988
$ws.Range("C9").Value = 1067.75
$ws.Range("D9").Value = -33.61002300565813

# SYSTEMS
$ws = $wb.Worksheets.Item("SYSTEMS")
$ws.Range("C2").Value = 1608.2999999999997
$ws.Range("D5").Value = 44.500404153453985
$ws.Range("D6").Value = 44.516335343245814
